$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Haver data pull refresh (6/24/21): the "date" code/reference row (row 2)
# is no longer part of the pulled series list, so remove it entirely and
# let every following row shift up one position.
$ws.Rows("2:2").Delete()

# Leave the selection where the editor left it after the refresh.
$ws.Range("G14").Select()
